$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Fix typo in the existing journal entry ("presquefini" -> "presque fini")
$ws.Range("B69").Value = "Travail en groupe pour faire fonctionner les budgets, les budgets partagés et les dettes. Quelques bugs restent à corriger, mais presque fini."

# Add the new journal entry in row 70
$ws.Range("A70").Value = 43238
$ws.Range("B70").Value = "Mis à jour de la branche fb-derby avec la branche master, on doit encore vérifier que tout marche avec Derby puis on merge les deux branches. On a essayer de créer le jar avec Guillaume, pas réussi. Problème de manifest, on check plus tard."
$ws.Range("C70").Value = 1.25

# The row grows taller to fit the wrapped text, like the other long entries above it
$ws.Rows.Item(70).RowHeight = 60

# Move the selection to where the user would naturally end up after typing the new row
[void]$ws.Range("B71").Select()
